# Update the "From" value for rule R30 (row 10, column C) on the Rules sheet
# from 18 to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("C10").Value = 1
